# Adds a "WINNING SCREEN" block (columns AB:AM) mirroring the existing
# INTRO/GAME screen blocks, with a centered "f\"PLAYER {player} WON!\""
# message box, per commit "Added check for winner functionality and winner screen".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlGeneral = 1

# ---------------------------------------------------------------------
# Row 2: "WINNING SCREEN" title, styled like the existing screen titles
# (C2:M2 "INTRO SCREEN", P2:Z2 "GAME SCREEN") - default font, centered.
# ---------------------------------------------------------------------
$title = $ws.Range("AC2:AM2")
$title.Merge()
$title.HorizontalAlignment = $xlCenter
$title.VerticalAlignment = $xlCenter
$ws.Range("AC2").Value = "WINNING SCREEN"

# ---------------------------------------------------------------------
# Row 3: thin divider/border row of 64s, like C3:M3 / P3:Z3
# (horizontal-center alignment only).
# ---------------------------------------------------------------------
$row3 = $ws.Range("AC3:AM3")
$row3.Value = 64
$row3.HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------
# Column AB: left divider marker column, like B and O (vertical-center
# alignment only).
# ---------------------------------------------------------------------
$colAB = $ws.Range("AB4:AB12")
$colAB.Value = 64
$colAB.VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------
# Row 4: thick top border of the screen box, like C4:M4.
# ---------------------------------------------------------------------
$ws.Range("AC4:AM4").BorderAround(1, -4138)
$ws.Range("AC4:AM4").Borders.Item(8).LineStyle = 1
$ws.Range("AC4:AM4").Borders.Item(8).Weight = -4138

# ---------------------------------------------------------------------
# Rows 5-11: plain box interior with thick left/right edges on AC / AM.
# ---------------------------------------------------------------------
$ws.Range("AC5:AC11").Borders.Item(7).LineStyle = 1
$ws.Range("AC5:AC11").Borders.Item(7).Weight = -4138
$ws.Range("AM5:AM11").Borders.Item(10).LineStyle = 1
$ws.Range("AM5:AM11").Borders.Item(10).Weight = -4138

# ---------------------------------------------------------------------
# Row 12: thick bottom border of the screen box, like C12:M12.
# ---------------------------------------------------------------------
$ws.Range("AC12:AM12").Borders.Item(9).LineStyle = 1
$ws.Range("AC12:AM12").Borders.Item(9).Weight = -4138
$ws.Range("AC12").Borders.Item(7).LineStyle = 1
$ws.Range("AC12").Borders.Item(7).Weight = -4138
$ws.Range("AM12").Borders.Item(10).LineStyle = 1
$ws.Range("AM12").Borders.Item(10).Weight = -4138

# ---------------------------------------------------------------------
# Inner "WON" message box: a plain vertically-centered ring around a
# merged AF7:AJ9 message cell (centered both ways).
# ---------------------------------------------------------------------
$ring = $ws.Range("AF6:AK6,AE7,AK7,AE8,AK8,AE9,AK9,AE10:AK10")
$ring.VerticalAlignment = $xlCenter

$msg = $ws.Range("AF7:AJ9")
$msg.Merge()
$msg.HorizontalAlignment = $xlCenter
$msg.VerticalAlignment = $xlCenter
$ws.Range("AF7").Value = 'f"PLAYER {player} WON!"'

# ---------------------------------------------------------------------
# View: scroll to show the new winning screen, select the message box
# header row, like the source file's AF5:AJ5 selection.
# ---------------------------------------------------------------------
$ws.Range("AF5:AJ5").Select()
